$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 82, shifting existing rows 82:223 down to 83:224
$ws.Rows(82).Insert()

# Populate the newly inserted row 82 with the new record's data
$ws.Cells.Item(82, 1).Value = 5
$ws.Cells.Item(82, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(82, 3).Value = "Maule"
$ws.Cells.Item(82, 4).Value = 44797
$ws.Cells.Item(82, 5).Value = 7
$ws.Cells.Item(82, 6).Value = 100112017
$ws.Cells.Item(82, 7).Value = "Apio"
$ws.Cells.Item(82, 8).Value = "Americana (o)"
$ws.Cells.Item(82, 9).Value = "Primera"
$ws.Cells.Item(82, 10).Value = 700
$ws.Cells.Item(82, 11).Value = 10000
$ws.Cells.Item(82, 12).Value = 10000
$ws.Cells.Item(82, 13).Value = 10000
$ws.Cells.Item(82, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(82, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(82, 16).Value = 833
$ws.Cells.Item(82, 17).Value = 12
$ws.Cells.Item(82, 18).Value = "Hortaliza"
